$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -149.2
$ws.Range("B3").Value = -138.3
$ws.Range("C3").Value = -109.4
$ws.Range("C4").Value = -104.4
$ws.Range("C10").Value = -164.3
$ws.Range("C21").Value = -150.9
$ws.Range("C23").Value = -4.7
